$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 67, pushing existing rows (old 67..154) down to 68..155
$ws.Rows.Item(67).Insert()

# Populate the newly inserted row 67 with the new data record
$ws.Cells.Item(67, 1).Value = 11
$ws.Cells.Item(67, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(67, 3).Value = 'Bíobío'
$ws.Cells.Item(67, 4).Value = 44705
$ws.Cells.Item(67, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(67, 5).Value = 8
$ws.Cells.Item(67, 6).Value = 'Fruta'
$ws.Cells.Item(67, 7).Value = 100109
$ws.Cells.Item(67, 8).Value = 'Uva'
$ws.Cells.Item(67, 9).Value = 100109001
$ws.Cells.Item(67, 10).Value = 'Uva'
$ws.Cells.Item(67, 11).Value = 'Red Globe'
$ws.Cells.Item(67, 12).Value = 'Primera'
$ws.Cells.Item(67, 13).Value = 220
$ws.Cells.Item(67, 14).Value = 9000
$ws.Cells.Item(67, 15).Value = 10000
$ws.Cells.Item(67, 16).Value = 9545
$ws.Cells.Item(67, 17).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(67, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(67, 19).Value = 530
$ws.Cells.Item(67, 20).Value = 18
